# Small improvements to thesis
# Add new feature rows (B19:B24) to the active worksheet's comparison list,
# and update the active cell selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "Communicate verbally"
$ws.Range("B20").Value = "Point to objects"
$ws.Range("B21").Value = "Independent of operating system and hardware platforms"
$ws.Range("B22").Value = "Multiple audio and video codecs integrated"
$ws.Range("B23").Value = "Integrated support for gestures"
$ws.Range("B24").Value = "Pause video feed"

[void]$ws.Range("B25").Select()
